# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# This updates the previously-computed "Strike#" values with the correct
# "K" (strikeout) counts for each outing, per the regenerated source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 4
    4  = 0
    5  = 4
    6  = 2
    7  = 5
    8  = 3
    9  = 0
    10 = 3
    11 = 0
    12 = 1
    13 = 3
    14 = 0
    15 = 3
    16 = 1
    17 = 1
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 2
    26 = 3
    27 = 0
    28 = 3
    29 = 0
    30 = 2
    31 = 2
    32 = 0
    33 = 0
    34 = 4
    35 = 4
    36 = 2
    37 = 0
    38 = 0
    39 = 1
    40 = 3
    41 = 2
    42 = 1
    43 = 3
    46 = 2
    47 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
